$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename header row strings:
#    *_old  -> *_FV2210
#    *_new  -> *_FV2304
# ---------------------------------------------------------------------------
$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------------
# 2. Turn the data range into a real table (ListObject) named "Table1".
#
#    The header row already carries custom formatting (bold / fill / border).
#    If that formatting is present while the table is created, Excel bakes it
#    into a header-row "dxf" override. To keep styles.xml untouched (as in
#    the target workbook) we stash the header formatting away, reset the
#    header to the default style, create the table, then restore the
#    formatting by copying it back (copying an existing style re-uses the
#    existing style record instead of synthesizing new ones).
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$stashRange = $ws.Range("A200:U200")

$headerRange.Copy()
$stashRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.Style = "Normal"

$tableRange = $ws.Range("A1:U86")
$listObject = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

$stashRange.Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$stashRange.Clear()

# ---------------------------------------------------------------------------
# 3. Freeze the header row (View > Freeze Panes > Freeze Top Row).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
